$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.238.83"
$ws.Range("E2").Value = "  +0.68%  "

$ws.Range("D3").Value = "3.743.12"
$ws.Range("E3").Value = "  +0.13%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "602.01"
$ws.Range("E5").Value = "  -0.05%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "168.01"
$ws.Range("E6").Value = "  +0.43%  "

$ws.Range("D7").Value = "3.739.90"
$ws.Range("E7").Value = "  +0.11%  "

$ws.Range("E8").Value = "  -0.05%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.542"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.170"
$ws.Range("E10").Value = "  +3.91%  "

$ws.Range("E11").Value = "  +0.14%  "

$ws.Range("E12").Value = "  +0.52%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "38.24"
$ws.Range("E13").Value = "  +0.61%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000249"
$ws.Range("E14").Value = "  +1.76%  "

$ws.Range("D15").Value = "4.365.52"
$ws.Range("E15").Value = "  +0.07%  "

$ws.Range("D16").Value = "3.735.75"
$ws.Range("E16").Value = "  -0.02%  "

$ws.Range("D17").Value = "69.057.72"
$ws.Range("E17").Value = "  +0.44%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.40"
$ws.Range("E18").Value = "  +2.13%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.38"
$ws.Range("E19").Value = "  +0.74%  "

$ws.Range("E20").Value = "  -1.43%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.23"
$ws.Range("E21").Value = "  +11.83%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "492.86"
$ws.Range("E22").Value = "  -0.76%  "

$ws.Range("E23").Value = "  +0.77%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000151"
$ws.Range("E24").Value = "  +6.10%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "84.86"
$ws.Range("E25").Value = "  -0.02%  "

$ws.Range("E26").Value = "  -0.53%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.31"
$ws.Range("E27").Value = "  -0.57%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.08"
$ws.Range("E28").Value = "  +0.18%  "

$ws.Range("E29").Value = "  -0.02%  "

$ws.Range("E30").Value = "  +1.53%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.16"
$ws.Range("E31").Value = "  +3.23%  "

$ws.Range("E32").Value = "  +1.16%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "31.68"
$ws.Range("E33").Value = "  -0.03%  "

$ws.Range("D34").Value = "3.884.14"
$ws.Range("E34").Value = "  +0.04%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.109"
$ws.Range("E35").Value = "  +0.72%  "

$ws.Range("D36").Value = "3.671.06"
$ws.Range("E36").Value = "  +0.02%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.998"
$ws.Range("E37").Value = "  -0.10%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.97"
$ws.Range("E38").Value = "  +2.63%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.140"
$ws.Range("E39").Value = "  +5.69%  "

$ws.Range("E40").Value = "  +0.02%  "

$ws.Range("E41").Value = "  +0.88%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.05"
$ws.Range("E42").Value = "  +6.32%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "48.84"
$ws.Range("E43").Value = "  -0.69%  "

$ws.Range("E44").Value = "  +1.14%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "424.84"
$ws.Range("E45").Value = "  -1.87%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.49"
$ws.Range("E46").Value = "  +0.89%  "

$ws.Range("E47").Value = "  +0.01%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "40.17"
$ws.Range("E48").Value = "  -1.16%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "141.19"
$ws.Range("E49").Value = "  +0.07%  "

$ws.Range("B50").Value = "VeChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0356"
$ws.Range("E50").Value = "  +0.85%  "

$ws.Range("B51").Value = "Maker"
$ws.Range("C51").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D51").Value = "2.784.79"
$ws.Range("E51").Value = "  +1.44%  "

